# AFDP-1080: Changing Brief title format
#
# The "Save Case File Rules" Drools decision table gets a new function
# (createTitle) used to build the Brief title from the case's Defendant
# last name + case number, plus the two new imports it needs
# (Person / PersonAssociation). The "Set Title as Case Number if null"
# rule is updated to call the new function, and its condition is widened
# so it always re-evaluates the title.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Make room for two new "Import" rows right after the existing
#    imports (before the "StandardEvaluationContext" import row), by
#    inserting two blank rows at row 13. Everything currently at row 13
#    and below (including the Functions block, the rule table, etc.)
#    shifts down by two rows and keeps its own formatting.
# ---------------------------------------------------------------------
$ws.Range("A13:A14").EntireRow.Insert()

$ws.Cells.Item(13, 3).Value2 = 'Import'
$ws.Cells.Item(13, 4).Value2 = 'com.armedia.acm.plugins.person.model.Person'

$ws.Cells.Item(14, 3).Value2 = 'Import'
$ws.Cells.Item(14, 4).Value2 = 'com.armedia.acm.plugins.person.model.PersonAssociation'

# ---------------------------------------------------------------------
# 2. Replace the Functions cell (now at row 16) with the combined
#    function block: the original dateFormat/addDays/evalSpring
#    functions (now as plain text, no more mixed rich-text run) plus
#    the new createTitle(CaseFile) helper.
# ---------------------------------------------------------------------
$functionsText = @'
    function String dateFormat(String fmt)
    {
      return LocalDate.now().toString(DateTimeFormat.forPattern(fmt));
    }
    function Date addDays(int days)
    {
      Calendar cal = Calendar.getInstance();
      cal.add(Calendar.DAY_OF_YEAR, days);
      return cal.getTime();
    }
    function Boolean evalSpring(String expression, Object obj)
    {
        ExpressionParser ep = new SpelExpressionParser();
        Expression exp = ep.parseExpression(expression);
        EvaluationContext ec = new StandardEvaluationContext();
        Boolean evaluated = exp.getValue(ec, obj, Boolean.class);
        return evaluated;
    }
    function String createTitle(CaseFile caseFile)
    {
        String lastName = "";
        String caseNumber = "";
        if (caseFile != null)
        {
            if (caseFile.getPersonAssociations() != null)
            {
                for (PersonAssociation pa : caseFile.getPersonAssociations())
                {
                   if ("Defendant".equalsIgnoreCase(pa.getPersonType()))
                   {
                      if (pa.getPerson() != null && pa.getPerson().getFamilyName() != null)
                      {
                         lastName = pa.getPerson().getFamilyName();
                         break;
                      }
                   }
               }
           }
           if (caseFile.getCaseNumber() != null)
           {
              caseNumber = caseFile.getCaseNumber();
           }
        }
        return lastName + "_" + caseNumber;
    }
'@

$ws.Cells.Item(16, 4).Value2 = $functionsText

# The cell now holds a lot more text, so give the row Excel's maximum
# row height (same as a user dragging it open / auto-fitting it).
$ws.Rows.Item(16).RowHeight = 409.5

# ---------------------------------------------------------------------
# 3. Update the "Set Title as Case Number if null" rule row (now at
#    row 29): widen the condition and switch the action to call the
#    new createTitle() function instead of the old literal prefix.
# ---------------------------------------------------------------------
$ws.Cells.Item(29, 3).Value2 = 'title == null || title != null'
$ws.Cells.Item(29, 4).Value2 = 'setTitle, createTitle($caseFile)'

# ---------------------------------------------------------------------
# 4. Refresh the view: scroll/selection like the author left it.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C29").Select()
